# ---------------------------------------------------------------------------
# Adds a "2022-Q4" sheet to the ASML workbook.
#
# Starting layout (tab order): 总计, 2022-Q3, 2022-Q2, 2022-Q1, 2020-Q4
# Target layout:                总计, 2022-Q4, 2022-Q3, 2022-Q2, 2022-Q1, 2020-Q4
#
# The former "2022-Q3" sheet is duplicated (so the new sheet inherits all of
# its formatting/column layout) and the duplicate becomes "2022-Q4" holding
# the new quarter's fund-holdings figures. The original "2022-Q3" sheet is
# left completely untouched (it simply slides one tab to the right). The
# "总计" (totals) roll-up sheet gets a new leading row for 2022-Q4 and all of
# its other rows shift down by one, with their counts/values updated to
# match.
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Helper: write a value into a cell as genuine TEXT (matching the workbook's
# convention of storing these figures as inline strings, not numbers) while
# leaving the cell's formatting/style untouched.
# ---------------------------------------------------------------------------
function Set-TextValue {
    param($range, [string]$text)
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.Style = "Normal"
}

# ---------------------------------------------------------------------------
# 1) Duplicate the current "2022-Q3" sheet (tab #2) to create the new
#    "2022-Q4" sheet just before it, then rename it.
# ---------------------------------------------------------------------------
$q3Sheet = $wb.Worksheets.Item(2)
$q3Sheet.Copy($q3Sheet, $null)

$q4Sheet = $wb.Worksheets.Item(2)
$q4Sheet.Name = "2022-Q4"

# ---------------------------------------------------------------------------
# 2) Overwrite the new "2022-Q4" sheet's fund rows with the Q4 figures.
# ---------------------------------------------------------------------------
Set-TextValue $q4Sheet.Range("B2") "270023"
Set-TextValue $q4Sheet.Range("C2") "广发全球精选股票（QDII）"
Set-TextValue $q4Sheet.Range("D2") "20.45"
Set-TextValue $q4Sheet.Range("E2") "82.63"
Set-TextValue $q4Sheet.Range("F2") "7.06"
Set-TextValue $q4Sheet.Range("G2") "1.4438"

Set-TextValue $q4Sheet.Range("B3") "000906"
Set-TextValue $q4Sheet.Range("C3") "广发全球精选股票（QDII）美元现汇"
Set-TextValue $q4Sheet.Range("D3") "20.45"
Set-TextValue $q4Sheet.Range("E3") "82.63"
Set-TextValue $q4Sheet.Range("F3") "7.06"
Set-TextValue $q4Sheet.Range("G3") "1.4438"

Set-TextValue $q4Sheet.Range("B4") "011423"
Set-TextValue $q4Sheet.Range("C4") "广发全球科技三个月定期开放混合（QDII）美元 C"
Set-TextValue $q4Sheet.Range("D4") "25.66"
Set-TextValue $q4Sheet.Range("E4") "89.07"
Set-TextValue $q4Sheet.Range("F4") "5.31"
Set-TextValue $q4Sheet.Range("G4") "1.3625"

Set-TextValue $q4Sheet.Range("D5") "21.02"
Set-TextValue $q4Sheet.Range("E5") "89.07"
Set-TextValue $q4Sheet.Range("F5") "5.31"
Set-TextValue $q4Sheet.Range("G5") "1.1162"

Set-TextValue $q4Sheet.Range("B6") "011422"
Set-TextValue $q4Sheet.Range("C6") "广发全球科技三个月定期开放混合（QDII）人民币 C"
Set-TextValue $q4Sheet.Range("D6") "4.84"
Set-TextValue $q4Sheet.Range("E6") "89.07"
Set-TextValue $q4Sheet.Range("F6") "5.31"
Set-TextValue $q4Sheet.Range("G6") "0.2570"

Set-TextValue $q4Sheet.Range("B7") "011421"
Set-TextValue $q4Sheet.Range("C7") "广发全球科技三个月定期开放混合（QDII）美元 A"
Set-TextValue $q4Sheet.Range("D7") "0.20"
Set-TextValue $q4Sheet.Range("E7") "89.07"
Set-TextValue $q4Sheet.Range("F7") "5.31"
Set-TextValue $q4Sheet.Range("G7") "0.0106"

# ---------------------------------------------------------------------------
# 3) Update the "总计" roll-up sheet: push every existing row down one slot
#    and fill in the new top row for 2022-Q4.
# ---------------------------------------------------------------------------
$totalSheet = $wb.Worksheets.Item(1)

# Give the new bottom row (A6:D6) the same look as the row above it before
# filling in its values.
$totalSheet.Range("A5").Copy()
$totalSheet.Range("A6").PasteSpecial(-4122)

$totalSheet.Range("A6").Value = 4
$totalSheet.Range("B6").Value = "2020-Q4"
$totalSheet.Range("C6").Value = 1
$totalSheet.Range("D6").Value = 0.02

$totalSheet.Range("B5").Value = "2022-Q1"
$totalSheet.Range("C5").Value = 5
$totalSheet.Range("D5").Value = 2.39

$totalSheet.Range("B4").Value = "2022-Q2"
$totalSheet.Range("C4").Value = 2
$totalSheet.Range("D4").Value = 1.37

$totalSheet.Range("B3").Value = "2022-Q3"
$totalSheet.Range("C3").Value = 6
$totalSheet.Range("D3").Value = 5.03

$totalSheet.Range("B2").Value = "2022-Q4"
$totalSheet.Range("C2").Value = 6
$totalSheet.Range("D2").Value = 5.63

# ---------------------------------------------------------------------------
# 4) Restore the originally-active tab ("2020-Q4", now the last sheet) since
#    duplicating the Q3 sheet shifted the active-sheet focus onto it.
# ---------------------------------------------------------------------------
$wb.Worksheets.Item($wb.Worksheets.Count).Activate()
